# added 4wk low sales check
$wb = $excel.ActiveWorkbook

# --- "Forecast Comparison" sheet: update MyForecast / Inventory Coverage / Seasonality Index ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$updates = @(
    @{ Row = 2;  D = 1; H = 19;   L = 1.18 },
    @{ Row = 3;  D = 1; H = 18;   L = 0.96 },
    @{ Row = 4;  D = 1; H = 17;   L = 1.02 },
    @{ Row = 5;  D = 1; H = 16;   L = 0.9399999999999999 },
    @{ Row = 6;  D = 1; H = 15;   L = 0.9399999999999999 },
    @{ Row = 7;  D = 1; H = 14;   L = 0.92 },
    @{ Row = 8;  D = 1; H = 13;   L = 0.8100000000000001 },
    @{ Row = 9;  D = 1; H = 12;   L = 0.96 },
    @{ Row = 10; D = 1; H = 11;   L = 1 },
    @{ Row = 11; D = 1; H = 10;   L = 1.11 },
    @{ Row = 12; D = 1; H = 9;    L = 0.96 },
    @{ Row = 13; D = 1; H = 8;    L = 1.09 },
    @{ Row = 14; D = 1; H = 7;    L = 0.86 },
    @{ Row = 15; D = 1; H = 6;    L = 0.87 },
    @{ Row = 16; D = 1; H = 5;    L = 0.9399999999999999 },
    @{ Row = 17; D = 1; H = 4;    L = 1.07 }
)

foreach ($u in $updates) {
    $ws1.Cells.Item($u.Row, 4).Value = $u.D   # column D = MyForecast
    $ws1.Cells.Item($u.Row, 8).Value = $u.H   # column H = Inventory Coverage
    $ws1.Cells.Item($u.Row, 12).Value = $u.L  # column L = Seasonality Index
}

# --- "Summary" sheet: update Total Forecast figures ---
# these cells hold numbers-as-text (e.g. "14"), so force a Text format
# before assigning, otherwise Excel would auto-coerce "16" etc. to a number
$ws2 = $wb.Worksheets.Item("Summary")

$ws2.Range("B9").NumberFormat = "@"
$ws2.Range("B9").Value = "16"

$ws2.Range("B10").NumberFormat = "@"
$ws2.Range("B10").Value = "8"

$ws2.Range("B11").NumberFormat = "@"
$ws2.Range("B11").Value = "4"
